# The "Diwali Offer 2022" (catId 4) and "Sunglasses" (catId 5) rows are cut
# from rows 5-6 and moved to the bottom of the table (rows 24-25); every row
# that was below them shifts up two rows to fill the gap. The "Catmaped
# with" value on the "Luxury Watch Collection" row (row 4) is corrected to
# point at itself instead of ".". Column B is widened and the active-cell
# selection is moved to B4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "Luxury Watch Collection"

$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "Watches"
$ws.Range("C5").Value = "https://cdn.cartpe.in/images/category_image_sm/5c598fb3a0d14.png"
$ws.Range("D5").Value = "https://watchcode1.cartpe.in/watches.html"

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = "Watches For Her"
$ws.Range("C6").Value = "https://cdn.cartpe.in/images/category_image_sm/659fed2dd8063.png"
$ws.Range("D6").Value = "https://watchcode1.cartpe.in/watches-for-her-watches.html"
$ws.Range("E6").Value = "Ladies Watch"

$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Perfume For Women"
$ws.Range("C7").Value = "https://cdn.cartpe.in/images/category_image_sm/659fea41962de.png"
$ws.Range("D7").Value = "https://watchcode1.cartpe.in/perfume-for-women-fragrance.html"

$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "Perfume For Men"
$ws.Range("C8").Value = "https://cdn.cartpe.in/images/category_image_sm/659fe9efea2a2.png"
$ws.Range("D8").Value = "https://watchcode1.cartpe.in/perfume-1.html"
$ws.Range("E8").Value = "."

$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Watches For Men"
$ws.Range("C9").Value = "https://cdn.cartpe.in/images/category_image_sm/659fed4129d32.png"
$ws.Range("D9").Value = "https://uzer-watch.cartpe.in/watches-for-men-watches.html"
$ws.Range("E9").Value = "Mens Watch"

$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "Ladies watch"
$ws.Range("C10").Value = "https://cdn.cartpe.in/images/category_image_sm/168992871764ba440d2c71d.jpeg"
$ws.Range("D10").Value = "https://city-watch19.cartpe.in/47703-ladies-watch.html"
$ws.Range("E10").Value = "Ladies Watch"

$ws.Range("A11").Value = 12
$ws.Range("B11").Value = "Girls Watch"
$ws.Range("C11").Value = "https://cartpe.in/images/no_image.jpg"
$ws.Range("D11").Value = "https://city-watch19.cartpe.in/girls-watch-women-accessories.html"
$ws.Range("E11").Value = "Ladies Watch"

$ws.Range("A12").Value = 13
$ws.Range("B12").Value = "Couple Watches"
$ws.Range("C12").Value = "https://cdn.cartpe.in/images/category_image_sm/5c5986f929154.jpg"
$ws.Range("D12").Value = "https://fashionpro123.cartpe.in/couple-watches.html"
$ws.Range("E12").Value = "."

$ws.Range("A13").Value = 14
$ws.Range("B13").Value = "G-Shock In SALE"
$ws.Range("D13").Value = "https://watch-aqua.cartpe.in/g-shock-in-sale-watches.html"
$ws.Range("E13").Value = "."

$ws.Range("A14").Value = 15
$ws.Range("B14").Value = "Wall Clock"
$ws.Range("C14").Value = "https://cdn.cartpe.in/images/category_image_sm/679ca9c0ead48.jpeg"
$ws.Range("D14").Value = "https://watchflex.cartpe.in/wall-clock-home-decor-home-home-.html"

$ws.Range("A15").Value = 16
$ws.Range("B15").Value = "Belts"
$ws.Range("C15").Value = "https://cdn.cartpe.in/images/category_image_sm/65c3e11bc47c5.jpeg"
$ws.Range("D15").Value = "https://mangoenterprise.cartpe.in/belts.html"

$ws.Range("A16").Value = 17
$ws.Range("B16").Value = "Wallet"
$ws.Range("C16").Value = "https://cdn.cartpe.in/images/category_image_sm/65c3e1243721e.jpeg"
$ws.Range("D16").Value = "https://mangoenterprise.cartpe.in/wallet.html"

$ws.Range("A17").Value = 18
$ws.Range("B17").Value = "Wallets and Belts"
$ws.Range("C17").Value = "https://cdn.cartpe.in/images/category_image_sm/65c3e12cd14ad.jpeg"
$ws.Range("D17").Value = "https://mangoenterprise.cartpe.in/wallets-amp-belts-men-accessories.html"

$ws.Range("A18").Value = 19
$ws.Range("B18").Value = "Ledish+Watch"
$ws.Range("C18").Value = "https://cdn.cartpe.in/images/category_image_sm/1684922045646ddebd3bbac.jpeg"
$ws.Range("D18").Value = "https://zeewatches.cartpe.in/.html"
$ws.Range("E18").Value = "Ladies Watch"

$ws.Range("A19").Value = 20
$ws.Range("B19").Value = "G-SHOCK+"
$ws.Range("C19").Value = "https://cdn.cartpe.in/images/category_image_sm/166209839063119bd664b36.jpeg"
$ws.Range("D19").Value = "https://eye-care.cartpe.in/43511-g-shock.html"

$ws.Range("A20").Value = 21
$ws.Range("B20").Value = "WALLET+BELT"
$ws.Range("C20").Value = "https://cdn.cartpe.in/images/category_image_sm/1662984721631f2211c6e14.jpeg"
$ws.Range("D20").Value = "https://eye-care.cartpe.in/45591-wallet.html"
$ws.Range("E20").Value = "."

$ws.Range("A21").Value = 22
$ws.Range("B21").Value = "Sunglasses and Frames"
$ws.Range("C21").Value = "https://cdn.cartpe.in/images/category_image_sm/63119bb0188e0.jpeg"
$ws.Range("D21").Value = "https://eye-care.cartpe.in/sunglasses-amp-frames-eye-wear-men.html"

$ws.Range("A22").Value = 23
$ws.Range("B22").Value = "Hand bags"
$ws.Range("C22").Value = "https://cdn.cartpe.in/images/category_image_sm/5c5987219e267.jpg"
$ws.Range("D22").Value = "https://eye-care.cartpe.in/hand-bags.html"

$ws.Range("A23").Value = 24
$ws.Range("B23").Value = "Home Decor"
$ws.Range("C23").Value = "https://cdn.cartpe.in/images/category_image_sm/679c5099e5cc7.jpeg"
$ws.Range("D23").Value = "https://eye-care.cartpe.in/home-decor-home-home-.html"

$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Diwali Offer 2022"
$ws.Range("C24").Value = "https://cartpe.in/images/no_image.jpg"
$ws.Range("D24").Value = "https://saenterprise.cartpe.in/diwali-offer-2022.html"

$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Sunglasses"
$ws.Range("C25").Value = "https://cdn.cartpe.in/images/category_image_sm/5c5995a7ce26e.jpeg"
$ws.Range("D25").Value = "https://thetimekeepers.cartpe.in/sunglasses-eye-wear-men.html"

# Column B ("catName") grows from 21.5703125 to 32.5703125 OOXML character-
# width units; the engine quantizes ColumnWidth input to the nearest 1/6 of
# a character, so 31.7 is the input that lands closest to the target width.
$ws.Columns("B").ColumnWidth = 31.7

# Active cell / selection moves from D31 to B4
$ws.Range("B4").Select()
